$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27
$ws.Range("A27").Value = "'5632"
$ws.Range("B27").Value = "'4/22/2025"
$ws.Range("C27").Value = "MORENO, JOSE MARIA AV. 93"
$ws.Range("D27").Value = "'6"
$ws.Range("E27").Value = "'804876047"
$ws.Range("H27").Value = "Aplomar"
$ws.Range("J27").Value = "Aplomo"
$ws.Range("L27").Value = "Pasante"
$ws.Range("M27").Value = -58.43607
$ws.Range("N27").Value = -34.61926

# Row 28
$ws.Range("A28").Value = "'5651"
$ws.Range("B28").Value = "'4/22/2025"
$ws.Range("C28").Value = "MONTES DE OCA, MANUEL AV. 511"
$ws.Range("D28").Value = "'4"
$ws.Range("E28").Value = "'804876051"
$ws.Range("H28").Value = "Picada"
$ws.Range("J28").Value = "Cambio"
$ws.Range("L28").Value = "Pasante"
$ws.Range("M28").Value = -58.375515
$ws.Range("N28").Value = -34.634393

# Row 29
$ws.Range("A29").Value = "'4562"
$ws.Range("B29").Value = "'4/23/2025"
$ws.Range("C29").Value = "LAS PALMAS 2620"
$ws.Range("D29").Value = "'4"
$ws.Range("E29").Value = "'804903802"
$ws.Range("H29").Value = "Poste inclinado"
$ws.Range("J29").Value = "Cambio"
$ws.Range("L29").Value = "Poste"
$ws.Range("M29").Value = -58.422686
$ws.Range("N29").Value = -34.647038

# Row 30
$ws.Range("A30").Value = "'5887"
$ws.Range("B30").Value = "'4/25/2025"
$ws.Range("C30").Value = "PALPA 3162"
$ws.Range("D30").Value = "'12"
$ws.Range("E30").Value = "'805010113"
$ws.Range("H30").Value = "Cambiar columna podrida en base."
$ws.Range("J30").Value = "Cambio"
$ws.Range("L30").Value = "Pasante"
$ws.Range("M30").Value = -58.451203
$ws.Range("N30").Value = -34.576561

# Row 31
$ws.Range("A31").Value = "'5671"
$ws.Range("B31").Value = "'4/28/2025"
$ws.Range("C31").Value = "TUCUMAN 3589"
$ws.Range("D31").Value = "'5"
$ws.Range("E31").Value = "'805507284"
$ws.Range("H31").Value = "Picada"
$ws.Range("J31").Value = "Cambio"
$ws.Range("L31").Value = "Pasante"
$ws.Range("M31").Value = -58.415839
$ws.Range("N31").Value = -34.599291

# Row 32
$ws.Range("A32").Value = "'5710"
$ws.Range("B32").Value = "'5/1/2025"
$ws.Range("C32").Value = "MONTES DE OCA, MANUEL AV. 141"
$ws.Range("D32").Value = "'4"
$ws.Range("E32").Value = "'805579077"
$ws.Range("H32").Value = "Picada"
$ws.Range("J32").Value = "Cambio"
$ws.Range("L32").Value = "Pasante"
$ws.Range("M32").Value = -58.376979
$ws.Range("N32").Value = -34.630568

# Row 33
$ws.Range("A33").Value = "'5716"
$ws.Range("B33").Value = "'5/1/2025"
$ws.Range("C33").Value = "NECOCHEA 1315"
$ws.Range("D33").Value = "'4"
$ws.Range("E33").Value = "'805579141"
$ws.Range("H33").Value = "Picada"
$ws.Range("J33").Value = "Cambio"
$ws.Range("L33").Value = "Pasante"
$ws.Range("M33").Value = -58.357009
$ws.Range("N33").Value = -34.635878

# Row 34
$ws.Range("A34").Value = "'5738"
$ws.Range("B34").Value = "'5/5/2025"
$ws.Range("C34").Value = "IRALA 29"
$ws.Range("D34").Value = "'4"
$ws.Range("E34").Value = "'805707145"
$ws.Range("H34").Value = "Picada e inclinada"
$ws.Range("J34").Value = "Cambio"
$ws.Range("L34").Value = "Terminal"
$ws.Range("M34").Value = -58.369244
$ws.Range("N34").Value = -34.628787

# Row 35
$ws.Range("A35").Value = "'5751"
$ws.Range("B35").Value = "'5/6/2025"
$ws.Range("C35").Value = "SALCEDO 2737"
$ws.Range("D35").Value = "'4"
$ws.Range("E35").Value = "'805707268"
$ws.Range("H35").Value = "Picada"
$ws.Range("J35").Value = "Cambio"
$ws.Range("L35").Value = "Pasante"
$ws.Range("M35").Value = -58.402664
$ws.Range("N35").Value = -34.631273

# Row 36
$ws.Range("A36").Value = "'6110"
$ws.Range("B36").Value = "'5/6/2025"
$ws.Range("C36").Value = "CORRALES 6147"
$ws.Range("D36").Value = "'8"
$ws.Range("E36").Value = "'805707291"
$ws.Range("H36").Value = "Columna corroida en su base"
$ws.Range("J36").Value = "Cambio"
$ws.Range("L36").Value = "Pasante"
$ws.Range("M36").Value = -58.469148
$ws.Range("N36").Value = -34.687883

# Row 37
$ws.Range("A37").Value = "'5823"
$ws.Range("B37").Value = "'5/19/2025"
$ws.Range("C37").Value = "CALVO, CARLOS AV. 3882"
$ws.Range("D37").Value = "'6"
$ws.Range("E37").Value = "'806926382"
$ws.Range("H37").Value = "picada"
$ws.Range("J37").Value = "Cambio"
$ws.Range("L37").Value = "Pasante"
$ws.Range("M37").Value = -58.419995
$ws.Range("N37").Value = -34.624709

# Row 38
$ws.Range("A38").Value = "'5838"
$ws.Range("B38").Value = "'5/19/2025"
$ws.Range("C38").Value = "ZAVALETA 105"
$ws.Range("D38").Value = "'4"
$ws.Range("E38").Value = "'806926388"
$ws.Range("H38").Value = "Aplomar"
$ws.Range("J38").Value = "Aplomo"
$ws.Range("L38").Value = "Pasante"
$ws.Range("M38").Value = -58.403556
$ws.Range("N38").Value = -34.638135

# Row 39
$ws.Range("A39").Value = "'5836"
$ws.Range("B39").Value = "'5/19/2025"
$ws.Range("C39").Value = "RIVADAVIA AV. 4548"
$ws.Range("D39").Value = "'6"
$ws.Range("E39").Value = "'806926405"
$ws.Range("H39").Value = "Aplomar"
$ws.Range("J39").Value = "Aplomo"
$ws.Range("L39").Value = "Pasante"
$ws.Range("M39").Value = -58.429977
$ws.Range("N39").Value = -34.615514

# Row 40
$ws.Range("A40").Value = "'5835"
$ws.Range("B40").Value = "'5/19/2025"
$ws.Range("C40").Value = "YAPEYU 198"
$ws.Range("D40").Value = "'5"
$ws.Range("E40").Value = "'806926444"
$ws.Range("H40").Value = "Aplomar"
$ws.Range("J40").Value = "Aplomo"
$ws.Range("L40").Value = "Pasante"
$ws.Range("M40").Value = -58.421623
$ws.Range("N40").Value = -34.614541

# Row 41
$ws.Range("A41").Value = "'-429"
$ws.Range("B41").Value = "'5/26/2025"
$ws.Range("C41").Value = "Blanco encalada 4362"
$ws.Range("D41").Value = "'12"
$ws.Range("E41").Value = "'806926710"
$ws.Range("H41").Value = "Cambiar columna 114 base corroida dar prioridad "
$ws.Range("J41").Value = "Cambio"
$ws.Range("L41").Value = "Pasante"
$ws.Range("M41").Value = -58.47888
$ws.Range("N41").Value = -34.571108

# Row 61: set OT (E61) which was previously empty
$ws.Range("E61").Value = "'807458227"
